# Update cryptos list values (price & 1h volume change) scraped on refresh.
# Rows 30/31 also swap Coin/Link (ranking order changed between ImmutableX and Kaspa).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '36.251.46'
$ws.Cells.Item(2, 5).Value = '  -3.69%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '1.970.76'
$ws.Cells.Item(3, 5).Value = '  -2.57%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
# Row 5
$ws.Cells.Item(5, 4).Value = '''230.06'
$ws.Cells.Item(5, 5).Value = '  -12.72%  '
# Row 6
$ws.Cells.Item(6, 4).Value = '''0.589'
$ws.Cells.Item(6, 5).Value = '  -4.95%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.06%  '
# Row 8
$ws.Cells.Item(8, 4).Value = '''53.69'
$ws.Cells.Item(8, 5).Value = '  -4.75%  '
# Row 9
$ws.Cells.Item(9, 4).Value = '''0.366'
$ws.Cells.Item(9, 5).Value = '  -5.30%  '
# Row 10
$ws.Cells.Item(10, 4).Value = '''57.17'
$ws.Cells.Item(10, 5).Value = '  +0.31%  '
# Row 11
$ws.Cells.Item(11, 4).Value = '''0.0744'
$ws.Cells.Item(11, 5).Value = '  -4.78%  '
# Row 12
$ws.Cells.Item(12, 4).Value = '''0.0980'
$ws.Cells.Item(12, 5).Value = '  -3.70%  '
# Row 13
$ws.Cells.Item(13, 4).Value = '2.259.96'
$ws.Cells.Item(13, 5).Value = '  -2.59%  '
# Row 14
$ws.Cells.Item(14, 4).Value = '''13.76'
$ws.Cells.Item(14, 5).Value = '  -4.72%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '''19.76'
$ws.Cells.Item(15, 5).Value = '  -5.25%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '''0.747'
$ws.Cells.Item(16, 5).Value = '  -7.73%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '''4.99'
$ws.Cells.Item(17, 5).Value = '  -5.36%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '1.965.21'
$ws.Cells.Item(18, 5).Value = '  -3.18%  '
# Row 19
$ws.Cells.Item(19, 4).Value = '36.196.81'
$ws.Cells.Item(19, 5).Value = '  -3.50%  '
# Row 20
$ws.Cells.Item(20, 4).Value = '''67.24'
$ws.Cells.Item(20, 5).Value = '  -3.80%  '
# Row 21
$ws.Cells.Item(21, 4).Value = '0.0₃0800'
$ws.Cells.Item(21, 5).Value = '  -5.44%  '
# Row 22
$ws.Cells.Item(22, 4).Value = '''5.05'
$ws.Cells.Item(22, 5).Value = '  -2.63%  '
# Row 23
$ws.Cells.Item(23, 4).Value = '''220.23'
$ws.Cells.Item(23, 5).Value = '  -3.58%  '
# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.02%  '
# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.83%  '
# Row 26
$ws.Cells.Item(26, 4).Value = '''2.32'
$ws.Cells.Item(26, 5).Value = '  -14.21%  '
# Row 27
$ws.Cells.Item(27, 4).Value = '''159.74'
$ws.Cells.Item(27, 5).Value = '  -2.67%  '
# Row 28
$ws.Cells.Item(28, 4).Value = '''8.47'
$ws.Cells.Item(28, 5).Value = '  -6.07%  '
# Row 29
$ws.Cells.Item(29, 4).Value = '''18.60'
$ws.Cells.Item(29, 5).Value = '  -5.57%  '
# Row 30
$ws.Cells.Item(30, 2).Value = 'Kaspa'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30, 4).Value = '''0.121'
$ws.Cells.Item(30, 5).Value = '  -5.60%  '
# Row 31
$ws.Cells.Item(31, 2).Value = 'ImmutableX'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31, 4).Value = '''1.31'
$ws.Cells.Item(31, 5).Value = '  -2.00%  '
# Row 32
$ws.Cells.Item(32, 5).Value = '  -3.88%  '
# Row 33
$ws.Cells.Item(33, 4).Value = '''4.31'
$ws.Cells.Item(33, 5).Value = '  -7.32%  '
# Row 34
$ws.Cells.Item(34, 4).Value = '''0.0599'
$ws.Cells.Item(34, 5).Value = '  -8.33%  '
# Row 35
$ws.Cells.Item(35, 4).Value = '''4.22'
$ws.Cells.Item(35, 5).Value = '  -7.36%  '
# Row 36
$ws.Cells.Item(36, 5).Value = '  -3.75%  '
# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.03%  '
# Row 38
$ws.Cells.Item(38, 4).Value = '''1.77'
$ws.Cells.Item(38, 5).Value = '  -2.99%  '
# Row 39
$ws.Cells.Item(39, 4).Value = '''3.19'
$ws.Cells.Item(39, 5).Value = '  -5.10%  '
# Row 40
$ws.Cells.Item(40, 4).Value = '''5.15'
$ws.Cells.Item(40, 5).Value = '  -1.76%  '
# Row 41
$ws.Cells.Item(41, 5).Value = '  -1.10%  '
# Row 42
$ws.Cells.Item(42, 4).Value = '1.425.78'
$ws.Cells.Item(42, 5).Value = '  +1.29%  '
# Row 43
$ws.Cells.Item(43, 4).Value = '''0.0878'
$ws.Cells.Item(43, 5).Value = '  -6.96%  '
# Row 44
$ws.Cells.Item(44, 4).Value = '''0.0199'
$ws.Cells.Item(44, 5).Value = '  -7.40%  '
# Row 45
$ws.Cells.Item(45, 5).Value = '  -12.55%  '
# Row 46
$ws.Cells.Item(46, 4).Value = '''87.16'
$ws.Cells.Item(46, 5).Value = '  -4.03%  '
# Row 47
$ws.Cells.Item(47, 4).Value = '''0.982'
$ws.Cells.Item(47, 5).Value = '  -4.94%  '
# Row 48
$ws.Cells.Item(48, 4).Value = '''14.63'
$ws.Cells.Item(48, 5).Value = '  -7.17%  '
# Row 49
$ws.Cells.Item(49, 4).Value = '''2.86'
$ws.Cells.Item(49, 5).Value = '  -0.60%  '
# Row 50
$ws.Cells.Item(50, 4).Value = '''6.72'
$ws.Cells.Item(50, 5).Value = '  -5.31%  '
# Row 51
$ws.Cells.Item(51, 4).Value = '''3.63'
$ws.Cells.Item(51, 5).Value = '  +14.21%  '
